$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fix typo in product name: "Nokia Lumio 1020" -> "Nokia Lumia 1020"
$ws.Range("A4").Value = "Nokia Lumia 1020"

# Move selection to A4 (last edited cell), matching the saved view state
$ws.Range("A4").Select()
